# Insert a new weekly price record as the new row 209 for
# "Hortaliza, Femacal de La Calera - Achicoria", pushing the previous
# rows 209:234 down to 210:235.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 209 (shifts 209:234 -> 210:235)
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new data record
$ws.Range("A209").Value = 3
$ws.Range("B209").Value = "Femacal de La Calera"
$ws.Range("C209").Value = "Coquimbo"
$ws.Range("D209").Value = 44858
$ws.Range("E209").Value = 5
$ws.Range("F209").Value = 100112010
$ws.Range("G209").Value = "Achicoria"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 60
$ws.Range("K209").Value = 5500
$ws.Range("L209").Value = 5500
$ws.Range("M209").Value = 5500
$ws.Range("N209").Value = "`$/caja 16 unidades"
$ws.Range("O209").Value = "Provincia de Quillota"
$ws.Range("P209").Value = 344
$ws.Range("Q209").Value = 16
$ws.Range("R209").Value = "Hortaliza"
